# Fix the "Sources" sheet content:
#  - correct the typo in the Employer Skills Survey description (A11):
#      "percent opf employers" -> "percent of employers"
#  - update the ESS 2023 rich-text note (B11) year from 2024 to 2023,
#      keeping "Coming Summer 2023" italicised like the original run
#  - move the active selection to A12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Sources")

# --- Update "Coming Summer 2024" -> "Coming Summer 2023" in the rich text, keeping italics ---
$b11 = $ws.Range("B11")
$b11Text = $b11.Value2
$yearIdx = $b11Text.IndexOf("2024")
if ($yearIdx -ge 0) {
    $yearChars = $b11.Characters($yearIdx + 1, 4)
    $yearChars.Text = "2023"
}
$b11Text2 = $b11.Value2
$comingIdx = $b11Text2.IndexOf("Coming")
if ($comingIdx -ge 0) {
    $comingChars = $b11.Characters($comingIdx + 1, $b11Text2.Length - $comingIdx)
    $comingChars.Font.Italic = $true
}

# --- Fix typo: "percent opf employers" -> "percent of employers" ---
$a11 = $ws.Range("A11")
$a11Text = $a11.Value2
$opfIdx = $a11Text.IndexOf("opf")
if ($opfIdx -ge 0) {
    $a11.Value2 = $a11Text.Substring(0, $opfIdx) + "of" + $a11Text.Substring($opfIdx + 3)
}

# --- Move selection to A12 ---
$ws.Range("A12").Select()
